# Regenerate save_data column G ("K") values: replace old Strike# counts with
# the recalculated K values (std/mean regen, s_vals calc/write).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 2
    4  = 0
    5  = 3
    6  = 2
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 2
    13 = 0
    14 = 1
    15 = 1
    16 = 3
    17 = 3
    18 = 3
    19 = 2
    20 = 5
    21 = 4
    22 = 2
    23 = 6
    24 = 11
    25 = 7
    26 = 6
    27 = 1
    28 = 2
    29 = 1
    30 = 2
    31 = 1
    32 = 2
    34 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
